$wb = $excel.ActiveWorkbook

# --- ELF-bldg-winter ---
$wsWinter = $wb.Worksheets.Item("ELF-bldg-winter")
$wsWinter.Range("B2").Value = 1.36302
$wsWinter.Range("D2").Value = 1.16033
$wsWinter.Range("B5").Value = 1.91535
$wsWinter.Range("D5").Value = 1.45616
$wsWinter.Range("D7").Value = 1.45616

# --- ELF-bldg-summer ---
$wsSummer = $wb.Worksheets.Item("ELF-bldg-summer")
$wsSummer.Range("B3").Value = 7.66676
$wsSummer.Range("D3").Value = 6.54006
$wsSummer.Range("B5").Value = 1.85042
$wsSummer.Range("D5").Value = 2.00709
$wsSummer.Range("D7").Value = 2.00709

# --- ELF-vehicles ---
$wsVehicles = $wb.Worksheets.Item("ELF-vehicles")
$wsVehicles.Range("B4").Value = 1.16038
$wsVehicles.Range("C4").Value = 1.22331
$wsVehicles.Range("B5").Value = 1.16038
$wsVehicles.Range("C5").Value = 1.22331
$wsVehicles.Range("B6").Value = 1.16038
$wsVehicles.Range("C6").Value = 1.22331
$wsVehicles.Range("B7").Value = 1.16038
$wsVehicles.Range("C7").Value = 1.22331

$wb.Save()
